$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting existing rows 29..72 down to 30..73.
$ws.Rows.Item(29).Insert(-4121)

# Populate the newly inserted row 29 with its data (mirrors the layout of the
# surrounding rows: Mercado ID, Mercado, Region, Fecha, Codreg, CategoriaID,
# Categoria, Variedad, Calidad, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Unidad, Origen, Precio $/Kg, Kg o Unidades, Clasificacion)
$ws.Range("A29").Value = 3
$ws.Range("B29").Value = "Femacal de La Calera"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = Get-Date -Year 2023 -Month 10 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Range("D29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = 300000000
$ws.Range("G29").Value = "Espárragos"
$ws.Range("H29").Value = "Verde"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 1300
$ws.Range("K29").Value = 1800
$ws.Range("L29").Value = 1800
$ws.Range("M29").Value = 1800
$ws.Range("N29").Value = "$/kilo"
$ws.Range("O29").Value = "Provincia de Quillota"
$ws.Range("P29").Value = 1800
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = "Hortaliza"
